$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to keep its value as literal text instead of letting
    # Excel auto-convert plain-number-looking strings into numeric values.
    $range.NumberFormat = "@"
    $range.Value = $value
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "65.898.22"
$ws.Range("E2").Value = "  +1.30%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.707.64"
$ws.Range("E3").Value = "  +2.81%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.04%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "609.93"
$ws.Range("E5").Value = "  +2.46%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "158.13"
$ws.Range("E6").Value = "  +1.36%  "

# Row 8 - XRP
Set-TextValue $ws.Range("D8") "0.589"
$ws.Range("E8").Value = "  -0.37%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  +5.71%  "

# Row 10 - Toncoin
$ws.Range("E10").Value = "  +3.83%  "

# Row 11 - Cardano
$ws.Range("E11").Value = "  +0.59%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  +1.07%  "

# Row 13 - Avalanche
Set-TextValue $ws.Range("D13") "30.45"
$ws.Range("E13").Value = "  +4.45%  "

# Row 14 - ShibaInu
Set-TextValue $ws.Range("D14") "0.0000204"
$ws.Range("E14").Value = "  +8.99%  "

# Row 15 - Wrapped liquid staked Ether 2.0
$ws.Range("D15").Value = "3.194.86"
$ws.Range("E15").Value = "  +2.70%  "

# Row 16 - Wrapped BTC
$ws.Range("D16").Value = "65.792.99"

# Row 17 - Wrapped Ether
$ws.Range("D17").Value = "2.704.08"
$ws.Range("E17").Value = "  +2.61%  "

# Row 18 - Chainlink
Set-TextValue $ws.Range("D18") "12.74"
$ws.Range("E18").Value = "  +1.52%  "

# Row 19 - Polkadot
$ws.Range("E19").Value = "  +2.18%  "

# Row 20 - Bitcoin Cash
Set-TextValue $ws.Range("D20") "360.09"
$ws.Range("E20").Value = "  +1.98%  "

# Row 21 - Uniswap
Set-TextValue $ws.Range("D21") "7.67"
$ws.Range("E21").Value = "  +3.93%  "

# Row 22 - Dai
$ws.Range("E22").Value = "  -0.08%  "

# Row 23 - Litecoin
Set-TextValue $ws.Range("D23") "70.95"
$ws.Range("E23").Value = "  +3.83%  "

# Row 24 - Internet Computer (DFINITY)
Set-TextValue $ws.Range("D24") "9.91"
$ws.Range("E24").Value = "  +3.98%  "

# Row 25 - PEPE
$ws.Range("E25").Value = "  +12.99%  "

# Row 26 - Sui Network
$ws.Range("E26").Value = "  -1.35%  "

# Row 27 - Fetch.AI
$ws.Range("E27").Value = "  +3.62%  "

# Row 28 - Kaspa
$ws.Range("E28").Value = "  +4.19%  "

# Row 29 - Aptos
Set-TextValue $ws.Range("D29") "8.42"
$ws.Range("E29").Value = "  +4.19%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  +5.27%  "

# Row 31 - Bittensor
Set-TextValue $ws.Range("D31") "545.87"
$ws.Range("E31").Value = "  +4.03%  "

# Row 32 - Binance-Peg BSC-USD
$ws.Range("E32").Value = "  +0.01%  "

# Row 33 - ImmutableX
Set-TextValue $ws.Range("D33") "1.83"
$ws.Range("E33").Value = "  +2.98%  "

# Row 34 - Render Token
Set-TextValue $ws.Range("D34") "6.72"
$ws.Range("E34").Value = "  +6.25%  "

# Row 35 - NEAR Protocol
Set-TextValue $ws.Range("D35") "5.46"
$ws.Range("E35").Value = "  -2.40%  "

# Row 36 - Polygon Ecosystem Token
$ws.Range("E36").Value = "  +2.22%  "

# Row 37 - Ethereum Classic
Set-TextValue $ws.Range("D37") "20.95"
$ws.Range("E37").Value = "  +3.28%  "

# Row 38 - Monero
Set-TextValue $ws.Range("D38") "163.53"
$ws.Range("E38").Value = "  +0.08%  "

# Row 39 - Stacks
$ws.Range("E39").Value = "  +0.00%  "

# Row 40 - First Digital USD
Set-TextValue $ws.Range("D40") "0.999"
$ws.Range("E40").Value = "  -0.02%  "

# Row 41 - Aave
Set-TextValue $ws.Range("D41") "173.59"
$ws.Range("E41").Value = "  +4.94%  "

# Row 43 - OKB
Set-TextValue $ws.Range("D43") "42.62"
$ws.Range("E43").Value = "  +0.85%  "

# Row 44 - Filecoin
$ws.Range("E44").Value = "  +3.08%  "

# Row 45 - Hedera
Set-TextValue $ws.Range("D45") "0.0618"
$ws.Range("E45").Value = "  +0.69%  "

# Row 46 - Injective Protocol
Set-TextValue $ws.Range("D46") "23.68"
$ws.Range("E46").Value = "  +2.52%  "

# Row 47 - dogwifhat
$ws.Range("E47").Value = "  +4.05%  "

# Row 48 - VeChain
$ws.Range("E48").Value = "  +4.39%  "

# Row 49 - Mantle
Set-TextValue $ws.Range("D49") "0.657"
$ws.Range("E49").Value = "  +1.54%  "

# Row 50 - EnergySwap
Set-TextValue $ws.Range("D50") "21.14"
$ws.Range("E50").Value = "  +8.77%  "

# Row 51 - Stellar
Set-TextValue $ws.Range("D51") "0.0994"
$ws.Range("E51").Value = "  +1.54%  "
